$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 0.1764705882352941
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 0.05263157894736842
$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 0.1111111111111111
$ws.Range("B8").Value = 1
$ws.Range("D8").Value = 0.05882352941176471
$ws.Range("B9").Value = 8
$ws.Range("D9").Value = 0.5333333333333333
$ws.Range("B10").Value = 3
$ws.Range("D10").Value = 0.1428571428571428
$ws.Range("B11").Value = 3
$ws.Range("D11").Value = 0.1578947368421053
$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 0.1052631578947368
$ws.Range("B13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("B14").Value = 6
$ws.Range("D14").Value = 0.3333333333333333
$ws.Range("B15").Value = 1
$ws.Range("D15").Value = 0.04761904761904762
$ws.Range("B16").Value = 1
$ws.Range("D16").Value = 0.05555555555555555
$ws.Range("B17").Value = 1
$ws.Range("D17").Value = 0.08333333333333333
$ws.Range("B19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("B22").Value = 1
$ws.Range("D22").Value = 0.07692307692307693
$ws.Range("B23").Value = 1
$ws.Range("D23").Value = 0.05
$ws.Range("B24").Value = 1
$ws.Range("D24").Value = 0.06666666666666667
$ws.Range("B25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("B27").Value = 1
$ws.Range("D27").Value = 0.06666666666666667
$ws.Range("B28").Value = 3
$ws.Range("D28").Value = 0.1764705882352941
$ws.Range("B29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("B32").Value = 8
$ws.Range("D32").Value = 0.4705882352941176
$ws.Range("B37").Value = 1
$ws.Range("D37").Value = 0.06666666666666667
$ws.Range("B39").Value = 6
$ws.Range("D39").Value = 0.2857142857142857
$ws.Range("B40").Value = 1
$ws.Range("D40").Value = 0.0625
$ws.Range("B41").Value = 2
$ws.Range("D41").Value = 0.1052631578947368
$ws.Range("B44").Value = 3
$ws.Range("D44").Value = 0.2
$ws.Range("B45").Value = 6
$ws.Range("D45").Value = 0.3529411764705883
$ws.Range("B46").Value = 2
$ws.Range("D46").Value = 0.1428571428571428
$ws.Range("B48").Value = 14
$ws.Range("D48").Value = 0.7
$ws.Range("B49").Value = 4
$ws.Range("D49").Value = 0.2857142857142857
$ws.Range("B50").Value = 2
$ws.Range("D50").Value = 0.125
$ws.Range("B51").Value = 4
$ws.Range("D51").Value = 0.2222222222222222
$ws.Range("B52").Value = 5
$ws.Range("D52").Value = 0.2777777777777778
$ws.Range("B54").Value = 9
$ws.Range("D54").Value = 0.5294117647058824
$ws.Range("B56").Value = 1
$ws.Range("D56").Value = 0.07142857142857142
$ws.Range("B57").Value = 22
$ws.Range("D57").Value = 1.222222222222222
$ws.Range("B58").Value = 8
$ws.Range("D58").Value = 0.6666666666666666
$ws.Range("B60").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("B61").Value = 1
$ws.Range("D61").Value = 0.0625
